$wb = $excel.ActiveWorkbook
$count = $wb.Worksheets.Count
$srcSheet = $wb.Worksheets.Item($count)
$newSheet = $wb.Worksheets.Add($null, $srcSheet)
$newSheet.Name = "2026-01-05"

$data = @(
    @('rank', 'title', 'author', 'latest_episode'),
    @('1', '貞操逆転世界で頼めばヤれると噂の俺', '澄田佑貴(漫画) aaa168（スリーエー）(原作)', '第4話'),
    @('2', 'いとこのこ', 'いぬちく(著者)', '第44話'),
    @('3', '新米オッサン冒険者、最強パーティに死ぬほど鍛えられて無敵になる', '漫画：荻野ケン 原作：岸馬きらく キャラクター原案：Tea', '【描きおろし公開】人気投票企画『みんなでハロウィン！』、結果発表！【第12巻発売記念】'),
    @('4', '異世界居酒屋「のぶ」', '蝉川夏哉(原作) ヴァージニア二等兵(漫画) 転(キャラクター原案)', '第128話'),
    @('5', '実は俺、最強でした？', '原作：澄守 彩 漫画：高橋 愛', 'おまけ68'),
    @('6', '田舎で恋は難しい!?', 'ねこうめ(著者)', '第3話'),
    @('7', '女友達は頼めば意外とヤらせてくれる', 'ろくろ(漫画) 鏡遊(原作)', '特別編'),
    @('8', '剥かせて！竜ケ崎さん', '一智和智', '大学生編 第16話'),
    @('9', '辺境モブ貴族のウチに嫁いできた悪役令嬢が、めちゃくちゃできる良い嫁なんだが？', 'tera(原作) 朝倉はやて(作画) 徹田(キャラクター原案)', 'お正月イラスト'),
    @('10', 'リビルドワールド', '綾村切人(漫画) ナフセ(原作) 吟(キャラクターデザイン) わいっしゅ(世界観デザイン) cell(メカニックデザイン)', '第76話④'),
    @('11', '小林さんちのメイドラゴン', 'クール教信者', '第157話'),
    @('12', '配信に致命的に向いていない女の子が迷宮で黙々と人助けする配信', '下田将也(漫画) 佐藤悪糖(原作) 福きつね(キャラクター原案)', '第5話後編'),
    @('13', '異世界食堂　洋食のねこや', '犬塚惇平(ヒーロー文庫／イマジカインフォス)(原作) ヤミザワ(漫画) モロザワ(漫画) エナミカツミ(キャラクター原案)', '第43話➁'),
    @('14', '願ってもない追放後からのスローライフ？ 〜引退したはずが成り行きで美少女ギャルの師匠になったらなぜかめちゃくちゃ懐かれた〜', 'ヤミーゴ(漫画) シュガースプーン。（GA文庫/SBクリエイティブ）(原作) なたーしゃ(キャラクター原案)', '第10話-1：友達と一緒なら'),
    @('15', 'アザミヤコを好きになる', 'ユニティコング(原作) ツノニガウ(作画)', '第12話後編'),
    @('16', '勇者パーティを追い出された器用貧乏　～パーティ事情で付与術士をやっていた剣士、万能へと至る～', '漫画：よねぞう 原作：都神樹 キャラクター原案：きさらぎゆり', '★描き下ろし色紙プレゼントのお知らせ★'),
    @('17', '理想のヒモ生活', '日月ネコ(漫画) 渡辺恒彦（ヒーロー文庫／イマジカインフォス）(原作) 文倉十(キャラクター原案)', '第91話　その1'),
    @('18', 'アイドル辞めるけど結婚してくれますか!?', '三吉汐美(著者)', '休載イラスト'),
    @('19', '落ちこぼれだった兄が実は最強 ～史上最強の勇者は転生し、学園で無自覚に無双する～', '村上よしゆき 茨木野 あるてら', '第４４話　勇者、S級ランクに認定される。弟は、ミカエルとクエストに向かったら、ラファエルが降ってくる（３）'),
    @('20', 'ウォルテニア戦記', '漫画：八木ゆかり 原作：保利亮太 キャラクター原案：bob', '第59話'),
    @('21', 'インフィニット・デンドログラム', '今井神 原作：海道左近 キャラクター原案：タイキ', '第75話'),
    @('22', '追放されたチート付与魔術師は 気ままなセカンドライフを謳歌する。', '六志麻あさ 業務用餅 kisui', '第７８話'),
    @('23', '十年目、帰還を諦めた転移者はいまさら主人公になる', '原作：氷純（「十年目、帰還を諦めた転移者はいまさら主人公になる」MFブックス刊） 漫画：しゅーかま キャラクター原案：あんべよしろう', '第21話③'),
    @('24', '姫様“拷問”の時間です', '原作:春原ロビンソン　漫画:ひらけい', '拷問164'),
    @('25', 'くらいあの子としたいこと', '碇マナツ(著者)', '年始特別編'),
    @('26', 'ずっと好きだった幼馴染と付き合い始めたら一途ビッチの性欲ジャンキーだったんだがどうすりゃいいですか？', '原作：トラ子猫 作画：あらいぐま', '第8話'),
    @('27', '黄金の経験値', '原純(原作) 霜月汐(作画) fixro2n(キャラクター原案)', '第21話（後編）'),
    @('28', '回復術士のやり直し', '月夜涙(原作) 羽賀ソウケン(漫画) しおこんぶ(キャラクター原案)', '第75話-2'),
    @('29', 'ダメ人間の愛しかた', '岩葉(著者)', '第23話後編　ダメ人間とプレゼントと彼女'),
    @('30', '最弱貴族に転生したので悪役たちを集めてみた', '空野進 sorani ファルまろ', '第15話　最弱貴族、領地を成長させる（４）'),
    @('31', 'モブ司祭だけど、この世界が乙女ゲームだと気づいたのでヒロインを育成します', '井冬良(漫画) レオナールＤ(原作) りいちゅ(キャラクター原案)', '第5話-1：祝福の時'),
    @('32', '全滅エンドを死に物狂いで回避した。パーティが病んだ。', '縞(漫画) 雨糸雀(原作) kodamazon(キャラクター原案)', '第2話後編：重すぎる想い'),
    @('33', 'ひとりぼっちの異世界攻略', 'びび（漫画） 五示正司（原作）', '第251話　エグい情報承りまーす'),
    @('34', '魔都精兵のスレイブ', '原作:タカヒロ　漫画:竹村洋平', '第169話　突入準備'),
    @('35', '理想の彼女', 'もりまりも(著者)', '第31話'),
    @('36', 'オタクに優しいギャルはいない!?', 'のりしろちゃん 魚住さかな', '【#172】サンタ予備軍'),
    @('37', '異世界帰りの剣聖は、自分の実力に気がつかない　～SSSランクの隠しクエストを受けて帰ってきたらレベル１のF級探索者になったけど、異世界で鍛え上げた剣術が強すぎて王女やS級英雄たちが俺を放ってくれません', '夜桜ユノ 吉田屋敷 もきゅ', '第２話　異世界帰りの剣聖は、ハイエナに「ざまぁ」する（３）'),
    @('38', '♀ガキとおじさん', 'サラマンダ(著者)', '第33話'),
    @('39', '生徒会役員共', '氏家ト全', '#429'),
    @('40', '異世界おじさん', '殆ど死んでいる(著者)', '第74話'),
    @('41', '35歳独身山田、異世界村に理想のセカンドハウスを作りたい　～異世界と現実のいいとこどりライフ～', '出雲大吉(原作) 西尾洋一(作画) ゆのひと(キャラクター原案)', '第6話②'),
    @('42', '銀河特急 ミルキー☆サブウェイ', '河野丼(漫画) 亀山陽平(原作・監修)', '第7話'),
    @('43', 'カナン様はあくまでチョロい', 'nonco', '第166話	カナンのチェストバースト'),
    @('44', '我輩は猫魔導師である ～キジトラ・ルークの快適チート猫生活～', '原作：猫神信仰研究会 漫画：三國大和 キャラクター原案：ハム', '第27話'),
    @('45', 'ノロマ魔法と呼ばれた魔法使いは重力魔法で無双する　～まだ重力の概念のない世界にて、少年は万有引力の王となる～', '神原絵理華(漫画) 一森一輝(原作)', '第22話③'),
    @('46', '一億年ボタンを連打した俺は、気付いたら最強になっていた ～落第剣士の学院無双～', '士土幽太郎(漫画) 月島秀一(原作) もきゅ(キャラクター原案)', '第43話-1'),
    @('47', '追放された底辺職「盗賊」はゲーム知識で無双する。一緒に召喚された先生も外れジョブだったけど効率的に成り上がります', 'ケンノジ 菅原イチバ カラスBTK', '第5話　温泉で無双する(2)'),
    @('48', '異種族レビュアーズ', '天原(原作) masha(作画)', '第89話'),
    @('49', 'バーサス', '原作：ONE 漫画：あずま京太郎 構成：bose', '★描き下ろし色紙プレゼントのお知らせ★'),
    @('50', '最強勇者パーティーは愛が知りたい', '山田肌襦袢', '第39話「デートは浮かれたっていい」')
)

for ($r = 0; $r -lt $data.Count; $r++) {
    for ($c = 0; $c -lt 4; $c++) {
        $newSheet.Cells.Item($r+1, $c+1).Value = $data[$r][$c]
    }
}

# Reuse the header-row formatting (bold, centered, bordered) from the
# previous weeks sheet so no new style entries are introduced.
$srcSheet.Range("A1:D1").Copy() | Out-Null
$newSheet.Range("A1:D1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

Write-Output ("Sheets: " + $wb.Worksheets.Count)
Write-Output ("New sheet name: " + $newSheet.Name)
